# feat: add 2022-Q1 data
#
# 1. Add a new "2022-Q1" worksheet (placed right before "总计", i.e. right
#    after "2021-Q4") with the quarter's fund-holding detail table.
# 2. Update the "总计" (totals) worksheet with a new leading row summarising
#    the 2022-Q1 quarter, shifting the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet.
#
#    We duplicate the existing "总计" sheet (rather than Worksheets.Add a
#    brand-new one) so that the page setup / outline properties / margins
#    that already exist on the workbook's sheets carry over faithfully, then
#    we drop its extra rows and overwrite the cell contents.
# ---------------------------------------------------------------------------
$totalWs0 = $wb.Worksheets.Item("总计")

# Places the duplicate immediately before "总计". NOTE: after Copy() runs,
# sheet references are resolved by live position, so $totalWs0 itself can no
# longer be trusted to mean the original "总计" sheet - fetch fresh handles
# for both sheets (by name) right after the copy completes.
$totalWs0.Copy($totalWs0, [System.Reflection.Missing]::Value)
$newWs = $wb.Worksheets.Item("总计 (2)")
$totalWs = $wb.Worksheets.Item("总计")
$newWs.Name = "2022-Q1"

# The old sheet had 5 data rows (rows 2-6); the new sheet only needs 3 (rows
# 2-4), so remove the two extra rows entirely.
$newWs.Range("A5:D6").EntireRow.Delete()

# Wipe the old values but keep the styling that was cloned along with the
# sheet (header/index-column formatting).
$newWs.Range("A1:D4").ClearContents()

# The old sheet only spanned columns A-D; extend the bold/bordered header
# styling out through column H to match the new, wider table.
$newWs.Range("D1").Copy()
$newWs.Range("E1:H1").PasteSpecial($xlPasteFormats)

# Header row
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Row 2
$newWs.Range("A2").Value = 0
$newWs.Range("B2").NumberFormat = "@"
$newWs.Range("B2").Value = "160222"
$newWs.Range("B2").Style = "Normal"
$newWs.Range("C2").Value = "国泰国证食品饮料行业指数（LOF）"
$newWs.Range("D2").NumberFormat = "@"
$newWs.Range("D2").Value = "57.34"
$newWs.Range("D2").Style = "Normal"
$newWs.Range("E2").NumberFormat = "@"
$newWs.Range("E2").Value = "92.33"
$newWs.Range("E2").Style = "Normal"
$newWs.Range("F2").NumberFormat = "@"
$newWs.Range("F2").Value = "1.62"
$newWs.Range("F2").Style = "Normal"
$newWs.Range("G2").NumberFormat = "@"
$newWs.Range("G2").Value = "0.9289"
$newWs.Range("G2").Style = "Normal"
$newWs.Range("H2").Value = 10

# Row 3
$newWs.Range("A3").Value = 1
$newWs.Range("B3").NumberFormat = "@"
$newWs.Range("B3").Value = "160323"
$newWs.Range("B3").Style = "Normal"
$newWs.Range("C3").Value = "华夏磐泰混合（LOF）"
$newWs.Range("D3").NumberFormat = "@"
$newWs.Range("D3").Value = "11.45"
$newWs.Range("D3").Style = "Normal"
$newWs.Range("E3").NumberFormat = "@"
$newWs.Range("E3").Value = "29.52"
$newWs.Range("E3").Style = "Normal"
$newWs.Range("F3").NumberFormat = "@"
$newWs.Range("F3").Value = "0.56"
$newWs.Range("F3").Style = "Normal"
$newWs.Range("G3").NumberFormat = "@"
$newWs.Range("G3").Value = "0.0641"
$newWs.Range("G3").Style = "Normal"
$newWs.Range("H3").Value = 4

# Row 4
$newWs.Range("A4").Value = 2
$newWs.Range("B4").NumberFormat = "@"
$newWs.Range("B4").Value = "161718"
$newWs.Range("B4").Style = "Normal"
$newWs.Range("C4").Value = "招商沪深300高贝塔指数"
$newWs.Range("D4").NumberFormat = "@"
$newWs.Range("D4").Value = "0.20"
$newWs.Range("D4").Style = "Normal"
$newWs.Range("E4").NumberFormat = "@"
$newWs.Range("E4").Value = "94.52"
$newWs.Range("E4").Style = "Normal"
$newWs.Range("F4").NumberFormat = "@"
$newWs.Range("F4").Value = "1.37"
$newWs.Range("F4").Style = "Normal"
$newWs.Range("G4").NumberFormat = "@"
$newWs.Range("G4").Value = "0.0027"
$newWs.Range("G4").Style = "Normal"
$newWs.Range("H4").Value = 6

# ---------------------------------------------------------------------------
# 2. Update "总计" sheet: insert a new top data row for "2022-Q1" and shift
#    the existing quarters' rows down by one, renumbering the index column.
# ---------------------------------------------------------------------------

# Clone the numbered index-column style onto the newly-needed row 7 before
# filling it in (rows 2-6 already own the correct style and simply get
# their values overwritten).
$totalWs.Range("A6").Copy()
$totalWs.Range("A7").PasteSpecial($xlPasteFormats)

# Fill bottom-up so every row's existing formatting stays in place while
# values shift down to make room for the new "2022-Q1" entry at the top.
$totalWs.Range("A7").Value = 5
$totalWs.Range("B7").Value = "2020-Q4"
$totalWs.Range("C7").Value = 29
$totalWs.Range("D7").Value = 9.619999999999999

$totalWs.Range("A6").Value = 4
$totalWs.Range("B6").Value = "2021-Q1"
$totalWs.Range("C6").Value = 6
$totalWs.Range("D6").Value = 1.85

$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q2"
$totalWs.Range("C5").Value = 11
$totalWs.Range("D5").Value = 2.7

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2021-Q3"
$totalWs.Range("C4").Value = 2
$totalWs.Range("D4").Value = 2.01

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2021-Q4"
$totalWs.Range("C3").Value = 7
$totalWs.Range("D3").Value = 1.58

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 3
$totalWs.Range("D2").Value = 1

# Restore the originally active sheet/tab ("2020-Q4") so we don't leave an
# unrelated side effect on sheetView/tabSelected state.
$wb.Worksheets.Item(1).Activate()
